$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix duplicated "类" in header text
$ws.Range("B1").Value = "其他服务类城市居民消费价格指数(上年=100)"
$ws.Range("D1").Value = "其他用品类城市居民消费价格指数(上年=100)"

# Copy the formatting of the last data row (A6) down onto the new rows
# so the new year labels pick up the same bold/centered/bordered style.
$ws.Range("A6").Copy()
$ws.Range("A7:A8").PasteSpecial(-4122)  # xlPasteFormats

# Add row 7: 2021年
$ws.Range("A7").Value = "2021年"
$ws.Range("B7").Value = 97.90000000000001
$ws.Range("C7").Value = 98.59999999999999
$ws.Range("D7").Value = 99.40000000000001

# Add row 8: 2022年 (B8 and D8 left blank/empty, but still present as text cells)
$ws.Range("A8").Value = "2022年"
$ws.Range("B8").Value = "'"
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").Value = 101.5
$ws.Range("D8").Value = "'"
$ws.Range("D8").Style = "Normal"
